$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A1 text from "source" to "text" (A2 keeps the long COVID report text)
$ws.Range("A1").Value = "text"

# Update selection to C2 (no frozen/top-left-cell override, i.e. default view)
[void]$ws.Range("C2").Select()
